$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab ("Worksheet" -> "Scores")
$ws.Name = "Scores"

# Column H only ever held empty, formatted placeholder cells - remove it
# entirely (this also renumbers the dimension/selection/spans and drops the
# now-unused "H" column width entry from <cols>).
$ws.Columns.Item(8).Delete()

# Re-apply column B's alignment/wrap formatting so the engine resolves it to
# the (now renumbered) matching style index instead of keeping the old one.
$dataRange = $ws.Range("B2:B169")
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.VerticalAlignment = -4108     # xlCenter
$dataRange.WrapText = $true

# Match the saved selection state from the target workbook
$ws.Range("I5").Select()
